# Apply the "Version correcciones sin escribir" edits to the metodologia workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Profesor / Dificultad ...): the cause description changes from
# "Dificultad para trabajar en grupo" to "Dificultad para conocer la
# aportación individual" (D11 "Carga inequitiva de trabjo" is unchanged).
$ws.Range("B11").Value = "Dificultad para conocer la aportación individual"

# New note in row 40 (column F) - added at the end of the data.
$ws.Range("F40").Value = "O como hacer hipotesis en usabilidad?"

# New note in B35 ("U otra"), the cell was previously blank (just styled).
$ws.Range("B35").Value = "U otra"

# New note in E11 ("Percepcion de carga nota"), previously blank cell.
$ws.Range("E11").Value = "Percepcion de carga nota"

# Update the view: scroll so row 8 is the top-left visible row, and the
# active selection becomes F11 instead of B30.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 8 | Out-Null
$win.ScrollColumn = 1 | Out-Null
$ws.Range("F11").Select() | Out-Null
